$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Add a new BOM row (row 31) for the shorting jumper part,
#    re-using the formatting of the previous last row (row 30).
# -----------------------------------------------------------------
$ws.Range("A30:G30").Copy()
$ws.Range("A31:G31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C31").Value = "SPC02SYAN"
$ws.Range("B31").Value = "CONN JUMPER SHORTING GOLD FLASH"
$ws.Range("A31").Value = "-"
$ws.Range("D31").Value = "http://www.sullinscorp.com/drawings/134_C02SYAN_11134.pdf"
$ws.Range("E31").Value = 0.083
$ws.Range("F31").Value = 9
$ws.Range("G31").Formula = "=F31*E31"

# Turn the datasheet reference into a real hyperlink, matching the
# rest of the "Datasheet" column, then restore the bordered hyperlink
# cell format that Hyperlinks.Add does not fully reproduce.
$ws.Hyperlinks.Add($ws.Range("D31"), "http://www.sullinscorp.com/drawings/134_C02SYAN_11134.pdf") | Out-Null
$ws.Range("D30").Copy()
$ws.Range("D31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# -----------------------------------------------------------------
# 2. Update the title in A1 (Power Board V1.2 -> Torpedo Board V3)
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Torpedo Board V3: Bill of Materials"

# -----------------------------------------------------------------
# 3. Misc view updates recorded the last time the sheet was saved.
# -----------------------------------------------------------------
$ws.Range("B19").Select()
$excel.ActiveWindow.Zoom = 100

Write-Host "Edit complete"
